# Atualizado por script em 05-11-2023 08:45
#
# denmark_3rd-division_2023-2024: the scraper re-ran and (a) two pairs of
# matches sharing the same kickoff date swapped their on-page order inside
# a handful of existing rows, and (b) four newly-played/newly-scraped
# fixtures were appended at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Re-order the match data (columns F:V) within six row-pairs -------
# Only the match columns move; Indice/pais/torneio/temporada/data_partida
# (A:E) stay attached to their original row position since both rows in
# each pair already share the same data_partida value.
$pairs = @(
    @(16, 17),
    @(21, 22),
    @(38, 39),
    @(41, 43),
    @(70, 71),
    @(73, 74)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $top = $ws.Range("F$($r1):V$($r1)").Value2
    $bottom = $ws.Range("F$($r2):V$($r2)").Value2
    $ws.Range("F$($r2):V$($r2)").Value = $top
    $ws.Range("F$($r1):V$($r1)").Value = $bottom
}

# --- 2) Append the four new fixtures as rows 80-83 ------------------------
# Mirror the formatting of the previous last row (79): bold/bordered index
# cell in A, and the datetime number format in E.
$ws.Range("A79").Copy() | Out-Null
$ws.Range("A80:A83").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E79").Copy() | Out-Null
$ws.Range("E80:E83").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Row 80: Lyseng 1-1 Naesby
$ws.Cells.Item(80, 1).Value = 79
$ws.Cells.Item(80, 2).Value = "denmark"
$ws.Cells.Item(80, 3).Value = "3rd-division"
$ws.Cells.Item(80, 4).Value = "2023-2024"
$ws.Cells.Item(80, 5).Value = 45234.54166666666
$ws.Cells.Item(80, 6).Value = "Lyseng"
$ws.Cells.Item(80, 7).Value = 1
$ws.Cells.Item(80, 8).Value = "Naesby"
$ws.Cells.Item(80, 9).Value = 1
$ws.Cells.Item(80, 10).Value = 2.09
$ws.Cells.Item(80, 11).Value = "03/11/2023 01:12"
$ws.Cells.Item(80, 12).Value = 2.41
$ws.Cells.Item(80, 13).Value = "04/11/2023 12:35"
$ws.Cells.Item(80, 14).Value = 3.44
$ws.Cells.Item(80, 15).Value = "03/11/2023 01:12"
$ws.Cells.Item(80, 16).Value = 3.2
$ws.Cells.Item(80, 17).Value = "04/11/2023 12:35"
$ws.Cells.Item(80, 18).Value = 2.78
$ws.Cells.Item(80, 19).Value = "03/11/2023 01:12"
$ws.Cells.Item(80, 20).Value = 2.8
$ws.Cells.Item(80, 21).Value = "04/11/2023 12:35"
$ws.Cells.Item(80, 22).Value = "https://www.betexplorer.com/football/denmark/3rd-division/lyseng-naesby/GOWlmiF1/"

# Row 81: Vanlose 0-3 SfB-Oure
$ws.Cells.Item(81, 1).Value = 80
$ws.Cells.Item(81, 2).Value = "denmark"
$ws.Cells.Item(81, 3).Value = "3rd-division"
$ws.Cells.Item(81, 4).Value = "2023-2024"
$ws.Cells.Item(81, 5).Value = 45234.54166666666
$ws.Cells.Item(81, 6).Value = "Vanlose"
$ws.Cells.Item(81, 7).Value = 0
$ws.Cells.Item(81, 8).Value = "SfB-Oure"
$ws.Cells.Item(81, 9).Value = 3
$ws.Cells.Item(81, 10).Value = 1.46
$ws.Cells.Item(81, 11).Value = "03/11/2023 01:12"
$ws.Cells.Item(81, 12).Value = 1.47
$ws.Cells.Item(81, 13).Value = "04/11/2023 12:59"
$ws.Cells.Item(81, 14).Value = 4.24
$ws.Cells.Item(81, 15).Value = "03/11/2023 01:12"
$ws.Cells.Item(81, 16).Value = 4.4
$ws.Cells.Item(81, 17).Value = "04/11/2023 12:59"
$ws.Cells.Item(81, 18).Value = 4.81
$ws.Cells.Item(81, 19).Value = "03/11/2023 01:12"
$ws.Cells.Item(81, 20).Value = 5.62
$ws.Cells.Item(81, 21).Value = "04/11/2023 12:59"
$ws.Cells.Item(81, 22).Value = "https://www.betexplorer.com/football/denmark/3rd-division/vanlose-sfb-oure/YczY9EME/"

# Row 82: Vejgaard 1-2 Young Boys
$ws.Cells.Item(82, 1).Value = 81
$ws.Cells.Item(82, 2).Value = "denmark"
$ws.Cells.Item(82, 3).Value = "3rd-division"
$ws.Cells.Item(82, 4).Value = "2023-2024"
$ws.Cells.Item(82, 5).Value = 45234.54166666666
$ws.Cells.Item(82, 6).Value = "Vejgaard"
$ws.Cells.Item(82, 7).Value = 1
$ws.Cells.Item(82, 8).Value = "Young Boys"
$ws.Cells.Item(82, 9).Value = 2
$ws.Cells.Item(82, 10).Value = 2.38
$ws.Cells.Item(82, 11).Value = "03/11/2023 01:12"
$ws.Cells.Item(82, 12).Value = 2.26
$ws.Cells.Item(82, 13).Value = "04/11/2023 12:43"
$ws.Cells.Item(82, 14).Value = 3.32
$ws.Cells.Item(82, 15).Value = "03/11/2023 01:12"
$ws.Cells.Item(82, 16).Value = 3.47
$ws.Cells.Item(82, 17).Value = "04/11/2023 12:43"
$ws.Cells.Item(82, 18).Value = 2.46
$ws.Cells.Item(82, 19).Value = "03/11/2023 01:12"
$ws.Cells.Item(82, 20).Value = 2.82
$ws.Cells.Item(82, 21).Value = "04/11/2023 12:43"
$ws.Cells.Item(82, 22).Value = "https://www.betexplorer.com/football/denmark/3rd-division/vejgaard-young-boys-fd/f7yUAf78/"

# Row 83: Holstebro 1-3 BK Frem
$ws.Cells.Item(83, 1).Value = 82
$ws.Cells.Item(83, 2).Value = "denmark"
$ws.Cells.Item(83, 3).Value = "3rd-division"
$ws.Cells.Item(83, 4).Value = "2023-2024"
$ws.Cells.Item(83, 5).Value = 45234.58333333334
$ws.Cells.Item(83, 6).Value = "Holstebro"
$ws.Cells.Item(83, 7).Value = 1
$ws.Cells.Item(83, 8).Value = "BK Frem"
$ws.Cells.Item(83, 9).Value = 3
$ws.Cells.Item(83, 10).Value = 3.04
$ws.Cells.Item(83, 11).Value = "03/11/2023 02:12"
$ws.Cells.Item(83, 12).Value = 3.01
$ws.Cells.Item(83, 13).Value = "03/11/2023 02:30"
$ws.Cells.Item(83, 14).Value = 3.34
$ws.Cells.Item(83, 15).Value = "03/11/2023 02:12"
$ws.Cells.Item(83, 16).Value = 3.4
$ws.Cells.Item(83, 17).Value = "04/11/2023 12:01"
$ws.Cells.Item(83, 18).Value = 2
$ws.Cells.Item(83, 19).Value = "03/11/2023 02:12"
$ws.Cells.Item(83, 20).Value = 2.16
$ws.Cells.Item(83, 21).Value = "03/11/2023 02:30"
$ws.Cells.Item(83, 22).Value = "https://www.betexplorer.com/football/denmark/3rd-division/holstebro-bk-frem/SKZy9YyL/"
